$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-15 20:44:56"
$wsZhCn.Range("K2").Value = "2016-08-15 20:45:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-15 20:45:05"
$wsDeDe.Range("K2").Value = "2016-08-15 20:45:37"
